$d = $word.ActiveDocument

# Remove the "Appendix: Quick prototype" appendix block:
#   - "Appendix: Quick prototype" heading paragraph
#   - the blank paragraph after it
#   - the "Figure: PDF page 1" paragraph
#   - the paragraph holding the embedded page-1.png drawing
# Locate the heading paragraph by its text and delete through the
# paragraph that contains the inline picture (the paragraph right
# before the next "Appendix: Links" Heading 2).

$startPara = $null
$endPara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "Appendix: Quick prototype") {
        $startPara = $i
        break
    }
}

if ($startPara -ne $null) {
    # Walk forward until (and including) the paragraph that contains the
    # inline picture; that is the last paragraph of the block to remove.
    for ($j = $startPara; $j -le $d.Paragraphs.Count; $j++) {
        $p = $d.Paragraphs.Item($j)
        if ($p.Range.InlineShapes.Count -gt 0) {
            $endPara = $j
            break
        }
    }

    if ($endPara -ne $null) {
        $rangeStart = $d.Paragraphs.Item($startPara).Range.Start
        $rangeEnd = $d.Paragraphs.Item($endPara).Range.End
        $r = $d.Range($rangeStart, $rangeEnd)
        $r.Delete()
    }
}
